# Generate Report for Handback
# The "7349fe32-fc3d-475b-a464-f06758fa10b0" entry now fails transform
# during handback, so it is resorted to directly follow "dfddac69..."
# (row 4) instead of sitting after "2d15fcfe..." (its previous slot,
# row 6). "15fb89bc..." and "2d15fcfe..." both shift down one row to
# make room, and the 7349fe32 row's status flips from
# "Ready for handoff" to "Handback transform failed".

$wb = $excel.ActiveWorkbook

function Set-CellAndHyperlink($ws, $addr, $value) {
    $ws.Range($addr).Value = $value
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $value
        }
    }
}

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
Set-CellAndHyperlink $ws '$A$4' "7349fe32-fc3d-475b-a464-f06758fa10b0.md"
$ws.Range("B4").Value = "Handback transform failed"
$ws.Range("C4").Value = "Handback transform failed"

Set-CellAndHyperlink $ws '$A$5' "15fb89bc-0cf8-4510-98da-b417cf8be285.md"
$ws.Range("B5").Value = "In Translation"
$ws.Range("C5").Value = "In Translation"

Set-CellAndHyperlink $ws '$A$6' "2d15fcfe-00a7-490d-9c8d-996bf8ae34bc.md"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
Set-CellAndHyperlink $ws '$A$4' "7349fe32-fc3d-475b-a464-f06758fa10b0.md"
$ws.Range("B4").Value = "Handback transform failed"
Set-CellAndHyperlink $ws '$C$4' "7349fe32-fc3d-475b-a464-f06758fa10b0.0099aaa95cfc73d6cae06941b67d6dec062ea2af.zh-cn.xlf"
$ws.Range("D4").Value = "2016-02-29 04:43:09"

Set-CellAndHyperlink $ws '$A$5' "15fb89bc-0cf8-4510-98da-b417cf8be285.md"
$ws.Range("B5").Value = "In Translation"
Set-CellAndHyperlink $ws '$C$5' "15fb89bc-0cf8-4510-98da-b417cf8be285.aa3d6acf6cc5302fcf6494129a9c8994df82172a.zh-cn.xlf"
$ws.Range("D5").Value = "2016-02-29 04:40:23"

Set-CellAndHyperlink $ws '$A$6' "2d15fcfe-00a7-490d-9c8d-996bf8ae34bc.md"
Set-CellAndHyperlink $ws '$C$6' "2d15fcfe-00a7-490d-9c8d-996bf8ae34bc.fa4766b55e040d881993370065a251c6581a92d4.zh-cn.xlf"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
Set-CellAndHyperlink $ws '$A$4' "7349fe32-fc3d-475b-a464-f06758fa10b0.md"
$ws.Range("B4").Value = "Handback transform failed"
Set-CellAndHyperlink $ws '$C$4' "7349fe32-fc3d-475b-a464-f06758fa10b0.0099aaa95cfc73d6cae06941b67d6dec062ea2af.de-de.xlf"
$ws.Range("D4").Value = "2016-02-29 04:43:24"

Set-CellAndHyperlink $ws '$A$5' "15fb89bc-0cf8-4510-98da-b417cf8be285.md"
$ws.Range("B5").Value = "In Translation"
Set-CellAndHyperlink $ws '$C$5' "15fb89bc-0cf8-4510-98da-b417cf8be285.aa3d6acf6cc5302fcf6494129a9c8994df82172a.de-de.xlf"
$ws.Range("D5").Value = "2016-02-29 04:40:34"

Set-CellAndHyperlink $ws '$A$6' "2d15fcfe-00a7-490d-9c8d-996bf8ae34bc.md"
Set-CellAndHyperlink $ws '$C$6' "2d15fcfe-00a7-490d-9c8d-996bf8ae34bc.fa4766b55e040d881993370065a251c6581a92d4.de-de.xlf"

$wb.Save()
